# Updates cryptos list cell values (Price/Volume columns) per latest data pull,
# and swaps the Chainlink / TRON rows (12 and 13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.919.78'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '2.301.77'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('D4').Formula = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Formula = "'300.15"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('D6').Formula = "'97.15"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.76%  '
$ws.Range('D7').Formula = "'0.512"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.86%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Formula = "'0.505"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.10%  '
$ws.Range('D10').Formula = "'35.59"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.49%  '
$ws.Range('D11').Formula = "'0.0787"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Formula = "'0.118"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Formula = "'17.90"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.20%  '
$ws.Range('D14').Formula = "'6.77"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.64%  '
$ws.Range('D15').Value = '2.659.50'
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').Value = '2.304.22'
$ws.Range('E16').Value = '  -5.67%  '
$ws.Range('D17').Formula = "'0.776"
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Value = '42.842.44'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').Formula = "'12.77"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.43%  '
$ws.Range('D20').Value = '0.0₃0902'
$ws.Range('E21').Value = '  -2.43%  '
$ws.Range('D22').Formula = "'67.85"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').Formula = "'240.02"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  -1.76%  '
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('D26').Formula = "'2.43"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Formula = "'25.36"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('D29').Formula = "'165.24"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.96%  '
$ws.Range('D30').Formula = "'2.02"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.35%  '
$ws.Range('D31').Formula = "'9.03"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.59%  '
$ws.Range('D32').Formula = "'32.94"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').Value = '  -3.60%  '
$ws.Range('D36').Formula = "'16.96"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.82%  '
$ws.Range('E37').Value = '  -1.21%  '
$ws.Range('D38').Formula = "'0.0685"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.50%  '
$ws.Range('E39').Value = '  -1.39%  '
$ws.Range('E40').Value = '  -2.89%  '
$ws.Range('E41').Value = '  -1.30%  '
$ws.Range('E42').Value = '  -2.04%  '
$ws.Range('D43').Value = '2.012.38'
$ws.Range('E43').Value = '  +0.75%  '
$ws.Range('E44').Value = '  -2.67%  '
$ws.Range('D45').Formula = "'10.15"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.44%  '
$ws.Range('D46').Formula = "'2.14"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('D47').Formula = "'17.45"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('E48').Value = '  -1.86%  '
$ws.Range('E49').Value = '  -2.23%  '
$ws.Range('D50').Value = '2.525.65'
$ws.Range('E50').Value = '  -0.55%  '
$ws.Range('D51').Formula = "'72.10"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.66%  '
